$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 208 (this shifts the existing rows 208-279
# down to 209-280, carrying formatting/styles with them, and updates the
# sheet dimension to A1:R280).
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row with the new price-record data.
$ws.Cells.Item(208, 1).Value = 4
$ws.Cells.Item(208, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(208, 3).Value = "Los Lagos"
$ws.Cells.Item(208, 4).Value = 44524
$ws.Cells.Item(208, 5).Value = 10
$ws.Cells.Item(208, 6).Value = 100114001
$ws.Cells.Item(208, 7).Value = "Papa"
$ws.Cells.Item(208, 8).Value = "Pehuenche"
$ws.Cells.Item(208, 9).Value = "1a nueva(o)"
$ws.Cells.Item(208, 10).Value = 150
$ws.Cells.Item(208, 11).Value = 13000
$ws.Cells.Item(208, 12).Value = 14000
$ws.Cells.Item(208, 13).Value = 13533
$ws.Cells.Item(208, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(208, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(208, 16).Value = 541
$ws.Cells.Item(208, 17).Value = 25
$ws.Cells.Item(208, 18).Value = "Hortaliza"
